$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the organization website value (was "www.stat.kg ", now "www.stat.gov.kg")
$ws.Range("B10").Value = "www.stat.gov.kg"

# Move the active selection to B10 (reflecting the edited cell)
$ws.Range("B10").Select()
